# Insert a new data row at row 99, pushing existing rows 99-141 down to 100-142,
# then populate the new row 99 with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 99:141 down by one row (inserting a new blank row 99).
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 44523
$ws.Range("D99").NumberFormat = $ws.Range("D100").NumberFormat
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 100114001
$ws.Range("G99").Value = "Papa"
$ws.Range("H99").Value = "Asterix"
$ws.Range("I99").Value = "1a nueva(o)"
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 11000
$ws.Range("M99").Value = 10500
$ws.Range("N99").Value = "`$/saco 25 kilos"
$ws.Range("O99").Value = "Provincia de Arauco"
$ws.Range("P99").Value = 420
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"
